$wb = $excel.ActiveWorkbook

# "BOM Report" sheet - Generated On timestamp
$wsBom = $wb.Worksheets.Item("BOM Report")
$wsBom.Range("B6").Value = "2020-02-09 9:09 PM"

# "Project Information" sheet - Report Time / Report Date / Report Date & Time
$wsInfo = $wb.Worksheets.Item("Project Information")
$wsInfo.Range("B8").Value = "9:09 PM"
$wsInfo.Range("B9").NumberFormat = "@"
$wsInfo.Range("B9").Value = "2020-02-09"
$wsInfo.Range("B10").Value = "2020-02-09 9:09 PM"
